$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1316
$ws.Range("I2").Value = 1354.3
$ws.Range("J2").Value = 1124.5
$ws.Range("K2").Value = 1354.3
$ws.Range("L2").Value = 1124.5
$ws.Range("M2").Value = -1241.3
$ws.Range("N2").Value = -1350.5
$ws.Range("H76").Value = 6617.364
$ws.Range("I76").Value = 5465.1665
$ws.Range("K76").Value = 5465.1665
$ws.Range("M76").Value = -5150.1665
$ws.Range("H79").Value = 6617.364
$ws.Range("I79").Value = 5465.1665
$ws.Range("K79").Value = 5465.1665
$ws.Range("M79").Value = -4373.1665
$ws.Range("H129").Value = 50002416
$ws.Range("I129").Value = 83334600
$ws.Range("K129").Value = 250003800
$ws.Range("M129").Value = -249998800
$ws.Range("H131").Value = 4368.591
$ws.Range("I131").Value = 2400.7144
$ws.Range("K131").Value = 7202.1432
$ws.Range("M131").Value = -2162.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3473496.2
$ws.Range("I2").Value = 4274472.5
$ws.Range("J2").Value = 2599.3333
$ws.Range("K2").Value = 4274472.5
$ws.Range("L2").Value = 2599.3333
$ws.Range("M2").Value = -4274359.5
$ws.Range("N2").Value = -2825.3333
$ws.Range("H50").Value = 2162.375
$ws.Range("I50").Value = 4024
$ws.Range("J50").Value = 1541.8334
$ws.Range("K50").Value = 4024
$ws.Range("L50").Value = 1541.8334
$ws.Range("M50").Value = -3310
$ws.Range("N50").Value = -2969.8334
$ws.Range("H116").Value = 3473496.2
$ws.Range("I116").Value = 4274472.5
$ws.Range("J116").Value = 2599.3333
$ws.Range("K116").Value = 4274472.5
$ws.Range("L116").Value = 2599.3333
$ws.Range("M116").Value = -4272178.5
$ws.Range("N116").Value = -7187.3333
$ws.Range("H122").Value = 580904
$ws.Range("I122").Value = 1716.4814
$ws.Range("J122").Value = 2318466.5
$ws.Range("K122").Value = 5149.4442
$ws.Range("L122").Value = 6955399.5
$ws.Range("M122").Value = -2699.4442
$ws.Range("N122").Value = -6960299.5
$ws.Range("H132").Value = 6777.759
$ws.Range("I132").Value = 7046.2856
$ws.Range("J132").Value = 6072.875
$ws.Range("K132").Value = 21138.8568
$ws.Range("L132").Value = 18218.625
$ws.Range("M132").Value = -18608.8568
$ws.Range("N132").Value = -23278.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3473496.2
$ws.Range("I3").Value = 4274472.5
$ws.Range("J3").Value = 2599.3333
$ws.Range("K3").Value = 4274472.5
$ws.Range("L3").Value = 2599.3333
$ws.Range("M3").Value = -4274358.5
$ws.Range("N3").Value = -2827.3333
$ws.Range("H86").Value = 5889237
$ws.Range("I86").Value = 20020880
$ws.Range("J86").Value = 1052.4166
$ws.Range("K86").Value = 20020880
$ws.Range("L86").Value = 1052.4166
$ws.Range("M86").Value = -20019757
$ws.Range("N86").Value = -3298.4166
$ws.Range("H89").Value = 5889237
$ws.Range("I89").Value = 20020880
$ws.Range("J89").Value = 1052.4166
$ws.Range("K89").Value = 100104400
$ws.Range("L89").Value = 5262.083000000001
$ws.Range("M89").Value = -100098784
$ws.Range("N89").Value = -16494.083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2652.5557
$ws.Range("I16").Value = 2268.2856
$ws.Range("K16").Value = 2268.2856
$ws.Range("M16").Value = -1981.2856
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 14629.769
$ws.Range("I31").Value = 6608.143
$ws.Range("K31").Value = 6608.143
$ws.Range("M31").Value = -6313.143
$ws.Range("H34").Value = 14629.769
$ws.Range("I34").Value = 6608.143
$ws.Range("K34").Value = 6608.143
$ws.Range("M34").Value = -6406.143
$ws.Range("H50").Value = 11916.667
$ws.Range("J50").Value = 12727.272
$ws.Range("L50").Value = 12727.272
$ws.Range("N50").Value = -13977.272
$ws.Range("H51").Value = 22999.334
$ws.Range("J51").Value = 29999
$ws.Range("L51").Value = 29999
$ws.Range("N51").Value = -31471
$ws.Range("H58").Value = 2599.4092
$ws.Range("I58").Value = 1788.0769
$ws.Range("K58").Value = 1788.0769
$ws.Range("M58").Value = -1585.0769
$ws.Range("H59").Value = 15000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 22999.334
$ws.Range("J61").Value = 29999
$ws.Range("L61").Value = 29999
$ws.Range("N61").Value = -30695
$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71498
$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -217488
$ws.Range("H86").Value = 13210
$ws.Range("I86").Value = 11596.429
$ws.Range("J86").Value = 14947.692
$ws.Range("K86").Value = 11596.429
$ws.Range("L86").Value = 14947.692
$ws.Range("M86").Value = -10473.429
$ws.Range("N86").Value = -17193.692
$ws.Range("H89").Value = 13210
$ws.Range("I89").Value = 11596.429
$ws.Range("J89").Value = 14947.692
$ws.Range("K89").Value = 57982.145
$ws.Range("L89").Value = 74738.45999999999
$ws.Range("M89").Value = -52366.145
$ws.Range("N89").Value = -85970.45999999999
$ws.Range("H113").Value = 2652.5557
$ws.Range("I113").Value = 2268.2856
$ws.Range("K113").Value = 2268.2856
$ws.Range("M113").Value = -98.28560000000016
$ws.Range("H136").Value = 2599.4092
$ws.Range("I136").Value = 1788.0769
$ws.Range("K136").Value = 5364.2307
$ws.Range("M136").Value = -2814.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 220
$ws.Range("I29").Value = 250
$ws.Range("J29").Value = 160
$ws.Range("K29").Value = 750
$ws.Range("L29").Value = 480
$ws.Range("M29").Value = -473
$ws.Range("N29").Value = -1034
$ws.Range("H92").Value = 1727.1818
$ws.Range("I92").Value = 2224.25
$ws.Range("J92").Value = 1443.1428
$ws.Range("K92").Value = 6672.75
$ws.Range("L92").Value = 4329.428400000001
$ws.Range("M92").Value = -5424.75
$ws.Range("N92").Value = -6825.428400000001
$ws.Range("H96").Value = 13686.923
$ws.Range("I96").Value = 9993.333000000001
$ws.Range("J96").Value = 14795
$ws.Range("K96").Value = 29979.999
$ws.Range("L96").Value = 44385
$ws.Range("M96").Value = -27920.999
$ws.Range("N96").Value = -48503
$ws.Range("H131").Value = 17363408
$ws.Range("I131").Value = 20833716
$ws.Range("J131").Value = 16669347
$ws.Range("K131").Value = 62501148
$ws.Range("L131").Value = 50008041
$ws.Range("M131").Value = -62496108
$ws.Range("N131").Value = -50018121

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H52").Value = 44606
$ws.Range("J52").Value = 44606
$ws.Range("L52").Value = 44606
$ws.Range("N52").Value = -45124
$ws.Range("H122").Value = 1784482.4
$ws.Range("I122").Value = 2971485.8
$ws.Range("K122").Value = 8914457.399999999
$ws.Range("M122").Value = -8912007.399999999
$ws.Range("H132").Value = 9256
$ws.Range("I132").Value = 6893.76
$ws.Range("J132").Value = 15817.777
$ws.Range("K132").Value = 20681.28
$ws.Range("L132").Value = 47453.331
$ws.Range("M132").Value = -18151.28
$ws.Range("N132").Value = -52513.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1465.5264
$ws.Range("J20").Value = 490
$ws.Range("L20").Value = 490
$ws.Range("N20").Value = -942
$ws.Range("H22").Value = 60358
$ws.Range("I22").Value = 124761.125
$ws.Range("J22").Value = 3110.7778
$ws.Range("K22").Value = 124761.125
$ws.Range("L22").Value = 3110.7778
$ws.Range("M22").Value = -124466.125
$ws.Range("N22").Value = -3700.7778
$ws.Range("H27").Value = 60358
$ws.Range("I27").Value = 124761.125
$ws.Range("J27").Value = 3110.7778
$ws.Range("K27").Value = 124761.125
$ws.Range("L27").Value = 3110.7778
$ws.Range("M27").Value = -124654.125
$ws.Range("N27").Value = -3324.7778
$ws.Range("H40").Value = 11499.2
$ws.Range("I40").Value = 10624.25
$ws.Range("K40").Value = 10624.25
$ws.Range("M40").Value = -10488.25
$ws.Range("H46").Value = 6404.0454
$ws.Range("I46").Value = 4732.8667
$ws.Range("K46").Value = 4732.8667
$ws.Range("M46").Value = -4544.8667
$ws.Range("H68").Value = 2860.5
$ws.Range("I68").Value = 3147.5
$ws.Range("J68").Value = 1999.5
$ws.Range("K68").Value = 3147.5
$ws.Range("L68").Value = 1999.5
$ws.Range("M68").Value = -2398.5
$ws.Range("N68").Value = -3497.5
$ws.Range("H71").Value = 2860.5
$ws.Range("I71").Value = 3147.5
$ws.Range("J71").Value = 1999.5
$ws.Range("K71").Value = 15737.5
$ws.Range("L71").Value = 9997.5
$ws.Range("M71").Value = -11993.5
$ws.Range("N71").Value = -17485.5
$ws.Range("H82").Value = 6174800
$ws.Range("I82").Value = 7938685.5
$ws.Range("J82").Value = 1200.5
$ws.Range("K82").Value = 7938685.5
$ws.Range("L82").Value = 1200.5
$ws.Range("M82").Value = -7938324.5
$ws.Range("N82").Value = -1922.5
$ws.Range("H85").Value = 6174800
$ws.Range("I85").Value = 7938685.5
$ws.Range("J85").Value = 1200.5
$ws.Range("K85").Value = 7938685.5
$ws.Range("L85").Value = 1200.5
$ws.Range("M85").Value = -7937437.5
$ws.Range("N85").Value = -3696.5
$ws.Range("H136").Value = 30825.783
$ws.Range("I136").Value = 42925.92
$ws.Range("J136").Value = 5617.1665
$ws.Range("K136").Value = 128777.76
$ws.Range("L136").Value = 16851.4995
$ws.Range("M136").Value = -126227.76
$ws.Range("N136").Value = -21951.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 17000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H61").Value = 11578.714
$ws.Range("I61").Value = 9410.4
$ws.Range("K61").Value = 9410.4
$ws.Range("M61").Value = -9118.4
